$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-19 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-20 Saturday", 2) | Out-Null
$d.Content.Find.Execute("230×3=690", $true, $false, $false, $false, $false, $true, 1, $false, "169×7=1183", 2) | Out-Null
$d.Content.Find.Execute("641×3=1923", $true, $false, $false, $false, $false, $true, 1, $false, "196×8=1568", 2) | Out-Null
$d.Content.Find.Execute("308×9=2772", $true, $false, $false, $false, $false, $true, 1, $false, "141×3=423", 2) | Out-Null
$d.Content.Find.Execute("291×9=2619", $true, $false, $false, $false, $false, $true, 1, $false, "659×6=3954", 2) | Out-Null
$d.Content.Find.Execute("510×2=1020", $true, $false, $false, $false, $false, $true, 1, $false, "850×5=4250", 2) | Out-Null
$d.Content.Find.Execute("364×6=2184", $true, $false, $false, $false, $false, $true, 1, $false, "804×9=7236", 2) | Out-Null
$d.Content.Find.Execute("131×2=262", $true, $false, $false, $false, $false, $true, 1, $false, "617×2=1234", 2) | Out-Null
$d.Content.Find.Execute("231×9=2079", $true, $false, $false, $false, $false, $true, 1, $false, "962×5=4810", 2) | Out-Null
$d.Content.Find.Execute("542×2=1084", $true, $false, $false, $false, $false, $true, 1, $false, "555×9=4995", 2) | Out-Null
$d.Content.Find.Execute("605×9=5445", $true, $false, $false, $false, $false, $true, 1, $false, "331×5=1655", 2) | Out-Null
$d.Content.Find.Execute("127×6=762", $true, $false, $false, $false, $false, $true, 1, $false, "831×3=2493", 2) | Out-Null
$d.Content.Find.Execute("472×5=2360", $true, $false, $false, $false, $false, $true, 1, $false, "712×4=2848", 2) | Out-Null
$d.Content.Find.Execute("770×9=6930", $true, $false, $false, $false, $false, $true, 1, $false, "461×3=1383", 2) | Out-Null
$d.Content.Find.Execute("413×5=2065", $true, $false, $false, $false, $false, $true, 1, $false, "705×8=5640", 2) | Out-Null
$d.Content.Find.Execute("693×5=3465", $true, $false, $false, $false, $false, $true, 1, $false, "439×8=3512", 2) | Out-Null
$d.Content.Find.Execute("593×8=4744", $true, $false, $false, $false, $false, $true, 1, $false, "270×6=1620", 2) | Out-Null
$d.Content.Find.Execute("495×2=990", $true, $false, $false, $false, $false, $true, 1, $false, "828×9=7452", 2) | Out-Null
$d.Content.Find.Execute("909×7=6363", $true, $false, $false, $false, $false, $true, 1, $false, "253×8=2024", 2) | Out-Null
$d.Content.Find.Execute("767×5=3835", $true, $false, $false, $false, $false, $true, 1, $false, "740×9=6660", 2) | Out-Null
$d.Content.Find.Execute("221×4=884", $true, $false, $false, $false, $false, $true, 1, $false, "844×2=1688", 2) | Out-Null
$d.Content.Find.Execute("942×2=1884", $true, $false, $false, $false, $false, $true, 1, $false, "706×7=4942", 2) | Out-Null
$d.Content.Find.Execute("857×3=2571", $true, $false, $false, $false, $false, $true, 1, $false, "836×5=4180", 2) | Out-Null
$d.Content.Find.Execute("461×5=2305", $true, $false, $false, $false, $false, $true, 1, $false, "556×4=2224", 2) | Out-Null
$d.Content.Find.Execute("113×8=904", $true, $false, $false, $false, $false, $true, 1, $false, "513×7=3591", 2) | Out-Null
$d.Content.Find.Execute("870×9=7830", $true, $false, $false, $false, $false, $true, 1, $false, "736×8=5888", 2) | Out-Null
